# "updated busted DK PN"
#
# Fix the wrong DigiKey part number for R5 (10k) on the main BOM sheet,
# and add a new "DK Order" sheet that only lists the Qty / Reference(s) /
# DK columns, for use when placing a DigiKey order.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The DigiKey part number "13-RC0402JR-1310KLTR-ND" was wrong (busted);
# replace it with the correct part number for R5 (10k).
$ws1.Range("D20").Value2 = "311-10KJRCT-ND"

# Add the new "DK Order" worksheet right after the main BOM sheet.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "DK Order"

# Column widths to roughly match the main sheet's Reference(s) / DK columns.
$ws2.Columns.Item(2).ColumnWidth = 12
$ws2.Columns.Item(3).ColumnWidth = 26

# Header row, bold like the main sheet's header.
$ws2.Range("A1").Value2 = "Qty"
$ws2.Range("B1").Value2 = "Reference(s)"
$ws2.Range("C1").Value2 = "DK"
$ws2.Range("A1:C1").Font.Bold = $true

$rows = @(
    @(3, "C1, C2, C3", "490-6539-1-ND"),
    @(1, "C4", "1276-6471-1-ND"),
    @(1, "C5", "1276-1552-1-ND"),
    @(1, "C6", "311-1114-1-ND"),
    @(1, "C7", "311-1128-1-ND"),
    @(1, "C8", "1276-1043-1-ND"),
    @(1, "C9", "1276-1096-1-ND"),
    @(1, "C10", "1292-1580-1-ND"),
    @(1, "C11", "P16200CT-ND"),
    @(1, "C12", "445-6008-1-ND"),
    @(4, "J1, J2, J3, J4", "A100886CT-ND"),
    @(1, "L1", "513-1568-1-ND"),
    @(1, "Q1", "296-25646-1-ND"),
    @(1, "Q2", "296-27625-1-ND"),
    @(1, "R1", "311-3.83KLRCT-ND"),
    @(1, "R2", "311-100KJRCT-ND"),
    @(1, "R3", "RMCF0402FT4K12CT-ND"),
    @(1, "R4", "311-1.33KLRCT-ND"),
    @(1, "R5", "311-10KJRCT-ND"),
    @(1, "R6", "RMCF0402FT5K62CT-ND"),
    @(1, "R7", "ST4ETB202CT-ND"),
    @(1, "U1", "296-25487-1-ND")
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value2 = $row[0]
    $ws2.Cells.Item($r, 2).Value2 = $row[1]
    $ws2.Cells.Item($r, 3).Value2 = $row[2]
    $r = $r + 1
}

# Restore selections / active sheet to match the saved workbook state.
[void]$ws2.Range("C29").Select()
[void]$ws1.Select()
[void]$ws1.Range("D27").Select()
